$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptocurrency price/volume update (GitHub Actions scheduled refresh)

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '42.948.99'
$ws.Range("E2").Value = '  +0.32%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.359.09'
$ws.Range("E3").Value = '  +1.54%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.01%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '302.91'
$ws.Range("E5").Value = '  +0.28%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '95.29'
$ws.Range("E6").Value = '  +1.40%  '

$ws.Range("E7").Value = '  +0.05%  '

$ws.Range("E8").Value = '  -0.32%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.477'
$ws.Range("E9").Value = '  -3.44%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '34.36'
$ws.Range("E10").Value = '  +1.23%  '

$ws.Range("E11").Value = '  +2.24%  '

$ws.Range("E12").Value = '  +0.48%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '18.43'
$ws.Range("E13").Value = '  -1.67%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.726.18'
$ws.Range("E14").Value = '  +1.53%  '

$ws.Range("E15").Value = '  -0.13%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.363.21'
$ws.Range("E16").Value = '  +1.76%  '

$ws.Range("E17").Value = '  +1.00%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '42.965.88'
$ws.Range("E18").Value = '  +0.54%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.90'
$ws.Range("E19").Value = '  -1.01%  '

$ws.Range("E20").Value = '  +1.60%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.0₃0884'
$ws.Range("E21").Value = '  -0.14%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '67.91'
$ws.Range("E22").Value = '  +0.06%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '235.13'
$ws.Range("E23").Value = '  -0.22%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.18'
$ws.Range("E24").Value = '  -1.96%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.44'
$ws.Range("E25").Value = '  +1.07%  '

$ws.Range("E26").Value = '  -0.06%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '24.37'
$ws.Range("E27").Value = '  -0.71%  '

$ws.Range("E28").Value = '  +15.57%  '

$ws.Range("E29").Value = '  +2.61%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '32.31'
$ws.Range("E30").Value = '  +3.03%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.00'
$ws.Range("E32").Value = '  +0.36%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '17.49'
$ws.Range("E33").Value = '  +0.36%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0724'
$ws.Range("E34").Value = '  +4.05%  '

$ws.Range("E35").Value = '  +6.30%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '127.85'
$ws.Range("E36").Value = '  -8.30%  '

$ws.Range("E37").Value = '  +0.65%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.32'
$ws.Range("E38").Value = '  -0.51%  '

$ws.Range("E39").Value = '  +3.54%  '

$ws.Range("E40").Value = '  -2.73%  '

$ws.Range("E41").Value = '  -0.71%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '20.89'
$ws.Range("E42").Value = '  -6.51%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.926.88'
$ws.Range("E43").Value = '  -0.33%  '

$ws.Range("E44").Value = '  +0.27%  '

$ws.Range("E45").Value = '  +3.38%  '

$ws.Range("E46").Value = '  -9.40%  '

$ws.Range("E47").Value = '  +0.38%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.587.38'
$ws.Range("E48").Value = '  +1.38%  '

$ws.Range("E49").Value = '  +2.94%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '71.36'
$ws.Range("E50").Value = '  -1.06%  '

$ws.Range("B51").Value = 'TrustWalletToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.13'
$ws.Range("E51").Value = '  +0.97%  '
